$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "328.00"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.58%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.99"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.22%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.250"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.21%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08087"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.14%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.522"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.04%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.654"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "5.07%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.05%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.13%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9355"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.37%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1332"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "24.76%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1971"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.47%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09085"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.34%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03472"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "5.04%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09584"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.02%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001352"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.86%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006402"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "7.63%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-6.78%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3514"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "3.04%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.003"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "12.78%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.67%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "7.44%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04449"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.18%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001223"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.63%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004286"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.08%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001202"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.04%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003995"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.08%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02483"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "12.73%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05182"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.74%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007650"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.60%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1424"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.54%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.009172"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "5.12%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002173"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.80%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01118"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "40.07%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006635"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.50%"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002483"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "148.08%"
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003343"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "16.80%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
